# Updates the crypto price/volume figures (and swaps the Fetch.AI / Dai
# rows) to match the latest GitHub Actions data refresh.
#
# Price-column values that look numeric are written with a leading
# apostrophe so Excel keeps them as text (matching the source data,
# which stores prices/percentages as plain strings, e.g. "6.60" rather
# than a numeric 6.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.818.50"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "2.910.54"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'586.79"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").Value = "'146.36"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.506"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "2.909.73"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").Value = "'6.97"
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").Value = "'0.435"
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "'32.80"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "3.391.68"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").Value = "61.845.93"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'6.60"
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("D19").Value = "2.908.39"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").Value = "'434.33"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").Value = "'80.93"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "'10.21"
$ws.Range("E26").Value = "  -7.91%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'2.08"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("E29").Value = "  +20.43%  "
$ws.Range("D30").Value = "'7.22"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'25.85"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Value = "'3.07"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").Value = "'5.50"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "'49.12"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").Value = "'38.86"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("D45").Value = "2.698.17"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'134.35"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'0.0337"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").Value = "'344.60"
$ws.Range("E48").Value = "  -8.13%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'22.33"
$ws.Range("E51").Value = "  -4.88%  "
